$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the expense table (header + 5 rows), adding the new
#     Item2/Qty2/Item3/Qty3 columns (F:I) and updating Item1 values ---

# Header row
$ws.Range("A1").Value = "Vendor"
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("D1").Value = "Item1"
$ws.Range("E1").Value = "Qty1"
$ws.Range("F1").Value = "Item2"
$ws.Range("G1").Value = "Qty2"
$ws.Range("H1").Value = "Item3"
$ws.Range("I1").Value = "Qty3"

# Row 2 - John Doe / Amazon
$ws.Range("A2").Value = "Amazon"
$ws.Range("B2").Value = "John "
$ws.Range("C2").Value = "Doe"
$ws.Range("D2").Value = "Watch"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = "Towel"
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = "Bike"
$ws.Range("I2").Value = 1

# Row 3 - Mary Smith / Walmart
$ws.Range("A3").Value = "Walmart"
$ws.Range("B3").Value = "Mary "
$ws.Range("C3").Value = "Smith"
$ws.Range("D3").Value = "Baseball"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "Shirt"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = "Pants"
$ws.Range("I3").Value = 1

# Row 4 - Tasha Locke / Target
$ws.Range("A4").Value = "Target"
$ws.Range("B4").Value = "Tasha"
$ws.Range("C4").Value = "Locke"
$ws.Range("D4").Value = "Candy"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = "Notebook"
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = "Camera"
$ws.Range("I4").Value = 1

# Row 5 - Hassan Baraka / Amazon
$ws.Range("A5").Value = "Amazon"
$ws.Range("B5").Value = "Hassan"
$ws.Range("C5").Value = "Baraka"
$ws.Range("D5").Value = "Watch"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "Bike"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = "Towel"
$ws.Range("I5").Value = 3

# Row 6 - Gabriel Alfaro / Walmart
$ws.Range("A6").Value = "Walmart"
$ws.Range("B6").Value = "Gabriel"
$ws.Range("C6").Value = "Alfaro"
$ws.Range("D6").Value = "Shirt"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = "Pants"
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = "Baseball"
$ws.Range("I6").Value = 3

# --- Column widths (manually resized by the author) ---
$ws.Columns.Item(1).ColumnWidth = 9.333
$ws.Columns.Item(2).ColumnWidth = 11.0
$ws.Columns.Item(3).ColumnWidth = 10.667
$ws.Columns.Item(6).ColumnWidth = 10.5
$ws.Columns.Item(7).ColumnWidth = 10.0
$ws.Columns.Item(8).ColumnWidth = 9.667

# --- Selection moved to D11 ---
$null = $ws.Range("D11").Select()
